$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.02105197657393327
$ws.Range("C2").Value = 0.1677375843530966

$ws.Range("B3").Value = 0.07045792989910823
$ws.Range("C3").Value = 0.2598345363540792

$ws.Range("B4").Value = 0.8609289576919302
$ws.Range("C4").Value = 0.4307109771549437

$ws.Range("B5").Value = 0.9924723883476512
$ws.Range("C5").Value = 0.3599945774992724

$ws.Range("B6").Value = 0.9720944538999806
$ws.Range("C6").Value = 0.7959176720437561

$ws.Range("B7").Value = 0.9340063778145274
$ws.Range("C7").Value = 0.2882401967297296

$ws.Range("B8").Value = 0.01279452443122864
$ws.Range("C8").Value = 0.1824441146850586
